$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row: "Activity Hour" / "theory"
$ws.Range("A136").Value = "Activity Hour"
$ws.Range("B136").Value = "theory"

# Fix the spelling of the existing "PROJ[TELAB/GTLAB]" entry (row 129, subject column)
$ws.Range("A129").Value = "PROJ[TE LAB/GT LAB]"

# New row: "Project" / "lab" - header-like row with a larger, explicit black font
$ws.Range("A137").Value = "Project"
$ws.Range("B137").Value = "lab"
$ws.Range("A137").Font.Size = 14
$ws.Range("A137").Font.Color = 0
$ws.Rows.Item(137).RowHeight = 18

# Update the view's active selection to match where the user's cursor ended up
$ws.Range("B138").Select()
